# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates column G ("K") values on Sheet1 for the 2023 Isaiah Campbell
# save-data workbook. The values are recomputed (strikeouts instead of the
# previous "Strike#" metric) and written back in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row number -> new value for column G (K)
$newKValues = [ordered]@{
    2  = 1
    3  = 2
    4  = 1
    5  = 2
    6  = 0
    7  = 1
    9  = 1
    10 = 3
    11 = 2
    12 = 1
    13 = 2
    14 = 1
    15 = 2
    16 = 0
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 1
    22 = 1
    23 = 2
    24 = 0
    25 = 0
    26 = 3
    27 = 2
    28 = 1
    29 = 1
}

foreach ($row in $newKValues.Keys) {
    $ws.Range("G$row").Value = $newKValues[$row]
}

$wb.Save()
